# Auto-generated edit applying scheduled market-data refresh to Aegis_Profits sheets.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per-row across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR worksheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 171966.5
$ws.Range("I40").Value = 509899.5
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 509899.5
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -509724.5
$ws.Range("N40").Value = -3350
$ws.Range("H43").Value = 3774.625
$ws.Range("I43").Value = 4533.3335
$ws.Range("J43").Value = 3319.4
$ws.Range("K43").Value = 4533.3335
$ws.Range("L43").Value = 3319.4
$ws.Range("M43").Value = -4464.3335
$ws.Range("N43").Value = -3457.4
$ws.Range("H52").Value = 732789.1
$ws.Range("I52").Value = 732789.1
$ws.Range("K52").Value = 2198367.3
$ws.Range("M52").Value = -2198207.3
$ws.Range("H58").Value = 1635857.8
$ws.Range("J58").Value = 2942.8572
$ws.Range("L58").Value = 8828.571599999999
$ws.Range("N58").Value = -9128.571599999999
$ws.Range("H64").Value = 74441.71000000001
$ws.Range("J64").Value = 3398
$ws.Range("L64").Value = 3398
$ws.Range("N64").Value = -3894
$ws.Range("H67").Value = 74441.71000000001
$ws.Range("J67").Value = 3398
$ws.Range("L67").Value = 3398
$ws.Range("N67").Value = -5114
$ws.Range("H74").Value = 6368.3335
$ws.Range("I74").Value = 7055.8335
$ws.Range("K74").Value = 7055.8335
$ws.Range("M74").Value = -6119.8335
$ws.Range("H77").Value = 6368.3335
$ws.Range("I77").Value = 7055.8335
$ws.Range("K77").Value = 35279.1675
$ws.Range("M77").Value = -30599.1675
$ws.Range("H98").Value = 623
$ws.Range("I98").Value = 513.375
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 513.375
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 984.625
$ws.Range("N98").Value = -4496
$ws.Range("H106").Value = 2253.4285
$ws.Range("I106").Value = 2016.9474
$ws.Range("K106").Value = 2016.9474
$ws.Range("M106").Value = -1385.9474
$ws.Range("H107").Value = 827.1111
$ws.Range("I107").Value = 816.94116
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 816.94116
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1103.05884
$ws.Range("N107").Value = -4840
$ws.Range("H113").Value = 60784.35
$ws.Range("I113").Value = 127016.75
$ws.Range("J113").Value = 1911.1111
$ws.Range("K113").Value = 127016.75
$ws.Range("L113").Value = 1911.1111
$ws.Range("M113").Value = -123762.75
$ws.Range("N113").Value = -8419.1111
$ws.Range("H118").Value = 8329.23
$ws.Range("I118").Value = 11436.667
$ws.Range("J118").Value = 1337.5
$ws.Range("K118").Value = 34310.001
$ws.Range("L118").Value = 4012.5
$ws.Range("M118").Value = -32653.001
$ws.Range("N118").Value = -7326.5
$ws.Range("H122").Value = 623
$ws.Range("I122").Value = 513.375
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 1540.125
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = 909.875
$ws.Range("N122").Value = -9400
$ws.Range("H125").Value = 1614.2
$ws.Range("I125").Value = 900
$ws.Range("J125").Value = 1792.75
$ws.Range("K125").Value = 8100
$ws.Range("L125").Value = 16134.75
$ws.Range("M125").Value = -5640
$ws.Range("N125").Value = -21054.75
$ws.Range("H129").Value = 2963.9375
$ws.Range("J129").Value = 909.2
$ws.Range("L129").Value = 2727.6
$ws.Range("N129").Value = -12727.6
$ws.Range("H137").Value = 1797.9166
$ws.Range("J137").Value = 2031.1818
$ws.Range("L137").Value = 6093.5454
$ws.Range("N137").Value = -11193.5454

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 12847.333
$ws.Range("I31").Value = 12847.333
$ws.Range("K31").Value = 12847.333
$ws.Range("M31").Value = -12553.333
$ws.Range("H32").Value = 9580.236999999999
$ws.Range("I32").Value = 7916.294
$ws.Range("J32").Value = 23723.75
$ws.Range("K32").Value = 7916.294
$ws.Range("L32").Value = 23723.75
$ws.Range("M32").Value = -7629.294
$ws.Range("N32").Value = -24297.75
$ws.Range("H45").Value = 54150.633
$ws.Range("J45").Value = 2356.3333
$ws.Range("L45").Value = 2356.3333
$ws.Range("N45").Value = -3110.3333
$ws.Range("H61").Value = 1260.5735
$ws.Range("I61").Value = 803.0294
$ws.Range("J61").Value = 1718.1177
$ws.Range("K61").Value = 803.0294
$ws.Range("L61").Value = 1718.1177
$ws.Range("M61").Value = -591.0294
$ws.Range("N61").Value = -2142.1177
$ws.Range("H74").Value = 1079.8286
$ws.Range("I74").Value = 1061
$ws.Range("J74").Value = 1120.909
$ws.Range("K74").Value = 1061
$ws.Range("L74").Value = 1120.909
$ws.Range("M74").Value = -187
$ws.Range("N74").Value = -2868.909
$ws.Range("H77").Value = 1079.8286
$ws.Range("I77").Value = 1061
$ws.Range("J77").Value = 1120.909
$ws.Range("K77").Value = 5305
$ws.Range("L77").Value = 5604.545
$ws.Range("M77").Value = -937
$ws.Range("N77").Value = -14340.545
$ws.Range("H136").Value = 1260.5735
$ws.Range("I136").Value = 803.0294
$ws.Range("J136").Value = 1718.1177
$ws.Range("K136").Value = 2409.0882
$ws.Range("L136").Value = 5154.3531
$ws.Range("M136").Value = 140.9117999999999
$ws.Range("N136").Value = -10254.3531

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 20778
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H134").Value = 3252.976
$ws.Range("I134").Value = 3292.0967
$ws.Range("J134").Value = 3142.7273
$ws.Range("K134").Value = 9876.2901
$ws.Range("L134").Value = 9428.1819
$ws.Range("M134").Value = -7341.2901
$ws.Range("N134").Value = -14498.1819

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2016.8334
$ws.Range("I134").Value = 933.1667
$ws.Range("K134").Value = 2799.5001
$ws.Range("M134").Value = -264.5001000000002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 231
$ws.Range("I10").Value = 119.75
$ws.Range("J10").Value = 676
$ws.Range("K10").Value = 359.25
$ws.Range("L10").Value = 2028
$ws.Range("M10").Value = -220.25
$ws.Range("N10").Value = -2306
$ws.Range("H47").Value = 174
$ws.Range("I47").Value = 123.333336
$ws.Range("J47").Value = 250
$ws.Range("K47").Value = 370.000008
$ws.Range("L47").Value = 750
$ws.Range("M47").Value = 60.99999200000002
$ws.Range("N47").Value = -1612
$ws.Range("H70").Value = 102681.2
$ws.Range("I70").Value = 501456
$ws.Range("J70").Value = 2987.5
$ws.Range("K70").Value = 1504368
$ws.Range("L70").Value = 8962.5
$ws.Range("M70").Value = -1504053
$ws.Range("N70").Value = -9592.5
$ws.Range("H73").Value = 102681.2
$ws.Range("I73").Value = 501456
$ws.Range("J73").Value = 2987.5
$ws.Range("K73").Value = 1504368
$ws.Range("L73").Value = 8962.5
$ws.Range("M73").Value = -1503276
$ws.Range("N73").Value = -11146.5
$ws.Range("H75").Value = 799.75
$ws.Range("J75").Value = 699.5
$ws.Range("L75").Value = 2098.5
$ws.Range("N75").Value = -4094.5
$ws.Range("H78").Value = 799.75
$ws.Range("J78").Value = 699.5
$ws.Range("L78").Value = 6295.5
$ws.Range("N78").Value = -16279.5
$ws.Range("H103").Value = 1366.0714
$ws.Range("I103").Value = 1082.5
$ws.Range("J103").Value = 1578.75
$ws.Range("K103").Value = 3247.5
$ws.Range("L103").Value = 4736.25
$ws.Range("M103").Value = -2368.5
$ws.Range("N103").Value = -6494.25
$ws.Range("H129").Value = 147989.08
$ws.Range("I129").Value = 6145.6
$ws.Range("J129").Value = 204726.48
$ws.Range("K129").Value = 18436.8
$ws.Range("L129").Value = 614179.4400000001
$ws.Range("M129").Value = -13436.8
$ws.Range("N129").Value = -624179.4400000001
$ws.Range("H131").Value = 824.51
$ws.Range("I131").Value = 431.34784
$ws.Range("J131").Value = 941.9480600000001
$ws.Range("K131").Value = 1294.04352
$ws.Range("L131").Value = 2825.84418
$ws.Range("M131").Value = 3745.95648
$ws.Range("N131").Value = -12905.84418
$ws.Range("H132").Value = 436578.8
$ws.Range("I132").Value = 810.4
$ws.Range("K132").Value = 7293.599999999999
$ws.Range("M132").Value = -4763.599999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3092.25
$ws.Range("I122").Value = 2455.5557
$ws.Range("J122").Value = 5002.3335
$ws.Range("K122").Value = 7366.6671
$ws.Range("L122").Value = 15007.0005
$ws.Range("M122").Value = -4916.6671
$ws.Range("N122").Value = -19907.0005
$ws.Range("H132").Value = 3996.25
$ws.Range("I132").Value = 3092.7273
$ws.Range("K132").Value = 9278.1819
$ws.Range("M132").Value = -6748.1819

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 25900
$ws.Range("I99").Value = 25900
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 25900
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -22905
$ws.Range("N99").ClearContents()
$ws.Range("H136").Value = 2744.7693
$ws.Range("I136").Value = 2055.3845
$ws.Range("J136").Value = 3434.1538
$ws.Range("K136").Value = 6166.1535
$ws.Range("L136").Value = 10302.4614
$ws.Range("M136").Value = -3616.1535
$ws.Range("N136").Value = -15402.4614

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 4625
$ws.Range("I132").Value = 2976.08
$ws.Range("J132").Value = 6498.773
$ws.Range("K132").Value = 8928.24
$ws.Range("L132").Value = 19496.319
$ws.Range("M132").Value = -6398.24
$ws.Range("N132").Value = -24556.319

